$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Add new test case row (TestCase_E42) to the bottom of the table ---
$newRow = 43

# Copy the formatting of the previous "paired-description" row (40) down to
# the new row first, so borders/wrap/row-height line up with the rest of the
# table, then overwrite the values.
$ws.Range("A40:E40").Copy()
$ws.Range("A" + $newRow + ":E" + $newRow).PasteSpecial(-4122)
$ws.Range("A" + $newRow).EntireRow.RowHeight = 30

# The description cell on this row uses the "shaded" description style
# (matching cells such as C5:C26), so pull that formatting in for column C.
$ws.Range("C5").Copy()
$ws.Range("C" + $newRow).PasteSpecial(-4122)

$ws.Range("A" + $newRow).Value = "TestCase_E42"
$ws.Range("C" + $newRow).Value = "Verify that user is able to watch an article to a particular watchlist from notification in home page||Verify that user is able to unwatch an article from watchlist from notification in home page"
$ws.Range("B" + $newRow).Value = "OPQA-298`n||OPQA-304"
$ws.Range("D" + $newRow).Value = "Y"
$ws.Range("E" + $newRow).Value = "PASS"

# --- Reset the view: scroll back to the top and select A2 ---
$ws.Range("A2").Select()
